$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "B1 + B2" label was renamed to "B1,B2" (same meaning, new format parsed
# by the updated parser per the commit message). It shows up in cells C3 and
# E3 - update both so the shared string is replaced (the engine drops the
# now-unused "B1 + B2" shared string automatically on save).
$ws.Range("C3").Value = "B1,B2"
$ws.Range("E3").Value = "B1,B2"

# Center-align the version header cells before merging, so the merged range
# keeps the centered style.
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("E1").HorizontalAlignment = -4108

# Merge the two version header cells (B1:C1 and D1:E1), matching the new
# two-column "Personal"/"Work" header layout.
$ws.Range("B1:C1").Merge()
$ws.Range("D1:E1").Merge()

# Move the active selection to E4 (was E8).
$null = $ws.Range("E4").Select()

Write-Host "Edit complete"
